# EPBDS-14346 Fix conversion when the argument type is String
#
# Adds a new "Method String str2str(String data)" test method (with body
# "return data;", mirroring the existing "oneArg" method) right after the
# "oneArg" method block and before the "Datatype Complex" table, on Sheet1
# of the multimodule test workbook.
#
# This is done by inserting 4 new rows after row 32 (pushing the existing
# "Datatype Complex" block from rows 33-36 down to rows 37-40), then
# populating the new rows 33-35 with the new method's header/body/blank-
# separator content (copying cell formatting from the analogous "oneArg"
# method block at rows 29-31), fixing up the block-start marker row (32)
# formatting, and finally re-creating the merged cells for the new method's
# header/body rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 33 (shifts old rows 33-36 -> 37-40).
$ws.Range("A33:E36").EntireRow.Insert()

# Merge the new method's header (B33:C33) and body (B34:C34) cells first;
# formats get pasted on top right after, so this doesn't disturb styling.
$ws.Range("B33:C33").Merge()
$ws.Range("B34:C34").Merge()

# Copy the "oneArg" method's header+body row formatting (rows 29-30) onto
# the new header+body rows (33-34).
$ws.Range("A29:E30").Copy()
$ws.Range("A33").PasteSpecial(-4122)

# Copy the blank separator row formatting (row 31) onto the new separator
# row (35).
$ws.Range("A31:E31").Copy()
$ws.Range("A35").PasteSpecial(-4122)

# Copy the blank row formatting (row 32) onto the trailing blank row (36).
$ws.Range("A32:E32").Copy()
$ws.Range("A36").PasteSpecial(-4122)

# Row 32 becomes the "start of block" marker row (same styling used before
# every other method block, e.g. B2:C2, B6:C6, B14:C14, ...).
$ws.Range("B2:C2").Copy()
$ws.Range("B32").PasteSpecial(-4122)

# Fill in the new method's text.
$ws.Range("B33").Value = "Method String str2str(String data)"
$ws.Range("B34").Value = "return data;"
